# "Ajout detection de mouvement fonctionnelle avec porte testés"
# The "alarme" row (row 4) on the Objectifs sheet is now functional/tested,
# so its "fait ?" status flips from "Non" to "Oui".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objectifs")

$ws.Activate()
$ws.Range("E4").Value = "Oui"

# Leave the selection where the author left it when they saved.
$ws.Range("G7").Select()
